$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.785.13"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.111.44"

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.88"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.55"

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.108.26"
$ws.Range("E8").Value = "  +0.94%  "

$ws.Range("E9").Value = "  -0.49%  "

$ws.Range("E10").Value = "  -1.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("E13").Value = "  -1.39%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.94"
$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.629.83"
$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.767.33"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.112.81"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.25"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "475.86"
$ws.Range("E21").Value = "  +2.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.93"
$ws.Range("E23").Value = "  +5.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.06"
$ws.Range("E24").Value = "  +0.98%  "

$ws.Range("E25").Value = "  +3.30%  "

$ws.Range("E26").Value = "  -2.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.84"
$ws.Range("E29").Value = "  -1.71%  "

$ws.Range("E30").Value = "  -1.54%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.53"
$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0937"
$ws.Range("E34").Value = "  -7.62%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.85"
$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.971"
$ws.Range("E37").Value = "  -3.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.04"

$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("E40").Value = "  -3.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.310"
$ws.Range("E41").Value = "  -2.29%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.67"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.831.23"
$ws.Range("E44").Value = "  +2.58%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "384.23"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("E46").Value = "  -1.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.57"
$ws.Range("E47").Value = "  -8.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.40"
$ws.Range("E48").Value = "  +0.21%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("E51").Value = "  -1.72%  "
